$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 208; existing rows 208..253 shift down to 209..254
$ws.Rows.Item(208).Insert()

# Populate the newly inserted row 208 with the new record's data.
# Columns A,B,C,E,F,G,H,I,J,R keep the same constant values used throughout
# this table; D,K,L,M,N,O,P,Q hold the new record's specific values.
$ws.Cells.Item(208, 1).Value = 4
$ws.Cells.Item(208, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(208, 3).Value = "Los Lagos"
$ws.Cells.Item(208, 4).Value = 44637
$ws.Cells.Item(208, 5).Value = 10
$ws.Cells.Item(208, 6).Value = 100112040
$ws.Cells.Item(208, 7).Value = "Cilantro"
$ws.Cells.Item(208, 8).Value = "Sin especificar"
$ws.Cells.Item(208, 9).Value = "Primera"
$ws.Cells.Item(208, 10).Value = 100
$ws.Cells.Item(208, 11).Value = 14000
$ws.Cells.Item(208, 12).Value = 14000
$ws.Cells.Item(208, 13).Value = 14000
$ws.Cells.Item(208, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(208, 15).Value = "Región Metropolitana"
$ws.Cells.Item(208, 16).Value = 389
$ws.Cells.Item(208, 17).Value = 36
$ws.Cells.Item(208, 18).Value = "Hortaliza"
